$wb = $excel.ActiveWorkbook

$survey   = $wb.Worksheets.Item("survey")
$settings = $wb.Worksheets.Item("settings")

# Rename the "display.title" setting label to "display.title.text"
$settings.Range("C1").Value = "display.title.text"

# Rename the "display.text" survey column header to "display.prompt.text"
$survey.Range("F1").Value = "display.prompt.text"

# Update the saved selections, finishing on "survey" so it becomes the
# active/selected tab (matching the workbook's new view state).
$settings.Range("C2").Select() | Out-Null
$survey.Range("F2").Select() | Out-Null
